# Refresh the auto-updating "datetimeFigureOut" date placeholders that live
# on the Slide Master, every Slide Layout, and the Notes Master so they show
# the current save date (10/8/18) instead of the stale 3/19/18 / 3/23/18
# values that were cached the last time the deck was saved.

$p = $ppt.ActivePresentation
$newDate = "10/8/18"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.Name -like "Date Placeholder*") {
            $shape.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide Master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every Slide Layout hanging off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Notes Master
Update-DatePlaceholder $p.NotesMaster.Shapes
